# Practice 5: remove the "Cliente" row from the Roles & Responsibilities table
# (ROL / RESPONSABILIDADES / TIPO DE CAMBIO / NIVEL DE ACCESO table).

$d = $word.ActiveDocument

$target = $null
foreach ($tbl in $d.Tables) {
    $headerText = $tbl.Cell(1, 1).Range.Text.TrimEnd([char]13, [char]7)
    if ($headerText -eq "ROL") {
        $target = $tbl
        break
    }
}

if ($target -ne $null) {
    $lastRow = $target.Rows.Count
    # Confirm the last row is the "Cliente" row before deleting it.
    $firstCellText = $target.Cell($lastRow, 1).Range.Text.TrimEnd([char]13, [char]7)
    if ($firstCellText -eq "Cliente") {
        $target.Rows($lastRow).Delete()
    }
}
